$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.174.55"
$ws.Range("E2").Value = "  -1.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.55"
$ws.Range("E3").Value = "  -1.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.38"
$ws.Range("E5").Value = "  -3.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4729"
$ws.Range("E7").Value = "  -2.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2748"
$ws.Range("E8").Value = "  -2.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06408"
$ws.Range("E9").Value = "  -1.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.860.38"
$ws.Range("E10").Value = "  -2.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07454"
$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("E12").Value = "  -2.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.997"
$ws.Range("E13").Value = "  -2.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.11"
$ws.Range("E14").Value = "  -4.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6291"
$ws.Range("E15").Value = "  -5.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.127.17"
$ws.Range("E16").Value = "  -1.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.63"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.70"
$ws.Range("E19").Value = "  -5.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007313"
$ws.Range("E20").Value = "  -4.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.095.94"
$ws.Range("E21").Value = "  -2.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.021"
$ws.Range("E23").Value = "  -5.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.987"
$ws.Range("E24").Value = "  -3.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.245"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.74"
$ws.Range("E26").Value = "  -1.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.80"
$ws.Range("E27").Value = "  -4.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.878"
$ws.Range("E28").Value = "  -3.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1018"
$ws.Range("E29").Value = "  +5.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.383"
$ws.Range("E30").Value = "  -2.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128"
$ws.Range("E31").Value = "  -5.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.915"
$ws.Range("E32").Value = "  -3.15%  "

$ws.Range("E33").Value = "  -3.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.142"
$ws.Range("E34").Value = "  -5.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7189"
$ws.Range("E35").Value = "  -4.22%  "

$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.692"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01888"
$ws.Range("E38").Value = "  +1.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.628"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9011"
$ws.Range("E40").Value = "  -1.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.952"
$ws.Range("E41").Value = "  -7.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.42"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4087"
$ws.Range("E44").Value = "  -4.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.524"
$ws.Range("E45").Value = "  -4.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.055"
$ws.Range("E46").Value = "  -5.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.05"
$ws.Range("E47").Value = "  -5.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1192"
$ws.Range("E48").Value = "  -7.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.681"
$ws.Range("E49").Value = "  -2.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.21"
$ws.Range("E50").Value = "  -2.07%  "

$ws.Range("E51").Value = "  -5.63%  "
